# Fruta / hortaliza, semanal
#
# The diff re-shuffles the per-row weekly data (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Origen,
# Precio $/Kg and Kg o Unidades - columns D, I, J, K, L, M, N, O, P, Q)
# across rows 3-16, leaving row 6 and all other columns untouched. It is a
# cyclic permutation of whole rows, so every original row is captured into
# variables first and only then written back out in its new position -
# this avoids any row being overwritten before its original value has been
# read.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for rows 3-16 (excluding row 6, which is unchanged)
# before permuting row contents (columns D,I,J,K,L,M,N,O,P,Q) according to the diff.
$D3 = $ws.Range("D3").Value()
$I3 = $ws.Range("I3").Value()
$J3 = $ws.Range("J3").Value()
$K3 = $ws.Range("K3").Value()
$L3 = $ws.Range("L3").Value()
$M3 = $ws.Range("M3").Value()
$N3 = $ws.Range("N3").Value()
$O3 = $ws.Range("O3").Value()
$P3 = $ws.Range("P3").Value()
$Q3 = $ws.Range("Q3").Value()

$D4 = $ws.Range("D4").Value()
$I4 = $ws.Range("I4").Value()
$J4 = $ws.Range("J4").Value()
$K4 = $ws.Range("K4").Value()
$L4 = $ws.Range("L4").Value()
$M4 = $ws.Range("M4").Value()
$N4 = $ws.Range("N4").Value()
$O4 = $ws.Range("O4").Value()
$P4 = $ws.Range("P4").Value()
$Q4 = $ws.Range("Q4").Value()

$D5 = $ws.Range("D5").Value()
$I5 = $ws.Range("I5").Value()
$J5 = $ws.Range("J5").Value()
$K5 = $ws.Range("K5").Value()
$L5 = $ws.Range("L5").Value()
$M5 = $ws.Range("M5").Value()
$N5 = $ws.Range("N5").Value()
$O5 = $ws.Range("O5").Value()
$P5 = $ws.Range("P5").Value()
$Q5 = $ws.Range("Q5").Value()

$D7 = $ws.Range("D7").Value()
$I7 = $ws.Range("I7").Value()
$J7 = $ws.Range("J7").Value()
$K7 = $ws.Range("K7").Value()
$L7 = $ws.Range("L7").Value()
$M7 = $ws.Range("M7").Value()
$N7 = $ws.Range("N7").Value()
$O7 = $ws.Range("O7").Value()
$P7 = $ws.Range("P7").Value()
$Q7 = $ws.Range("Q7").Value()

$D8 = $ws.Range("D8").Value()
$I8 = $ws.Range("I8").Value()
$J8 = $ws.Range("J8").Value()
$K8 = $ws.Range("K8").Value()
$L8 = $ws.Range("L8").Value()
$M8 = $ws.Range("M8").Value()
$N8 = $ws.Range("N8").Value()
$O8 = $ws.Range("O8").Value()
$P8 = $ws.Range("P8").Value()
$Q8 = $ws.Range("Q8").Value()

$D9 = $ws.Range("D9").Value()
$I9 = $ws.Range("I9").Value()
$J9 = $ws.Range("J9").Value()
$K9 = $ws.Range("K9").Value()
$L9 = $ws.Range("L9").Value()
$M9 = $ws.Range("M9").Value()
$N9 = $ws.Range("N9").Value()
$O9 = $ws.Range("O9").Value()
$P9 = $ws.Range("P9").Value()
$Q9 = $ws.Range("Q9").Value()

$D10 = $ws.Range("D10").Value()
$I10 = $ws.Range("I10").Value()
$J10 = $ws.Range("J10").Value()
$K10 = $ws.Range("K10").Value()
$L10 = $ws.Range("L10").Value()
$M10 = $ws.Range("M10").Value()
$N10 = $ws.Range("N10").Value()
$O10 = $ws.Range("O10").Value()
$P10 = $ws.Range("P10").Value()
$Q10 = $ws.Range("Q10").Value()

$D11 = $ws.Range("D11").Value()
$I11 = $ws.Range("I11").Value()
$J11 = $ws.Range("J11").Value()
$K11 = $ws.Range("K11").Value()
$L11 = $ws.Range("L11").Value()
$M11 = $ws.Range("M11").Value()
$N11 = $ws.Range("N11").Value()
$O11 = $ws.Range("O11").Value()
$P11 = $ws.Range("P11").Value()
$Q11 = $ws.Range("Q11").Value()

$D12 = $ws.Range("D12").Value()
$I12 = $ws.Range("I12").Value()
$J12 = $ws.Range("J12").Value()
$K12 = $ws.Range("K12").Value()
$L12 = $ws.Range("L12").Value()
$M12 = $ws.Range("M12").Value()
$N12 = $ws.Range("N12").Value()
$O12 = $ws.Range("O12").Value()
$P12 = $ws.Range("P12").Value()
$Q12 = $ws.Range("Q12").Value()

$D13 = $ws.Range("D13").Value()
$I13 = $ws.Range("I13").Value()
$J13 = $ws.Range("J13").Value()
$K13 = $ws.Range("K13").Value()
$L13 = $ws.Range("L13").Value()
$M13 = $ws.Range("M13").Value()
$N13 = $ws.Range("N13").Value()
$O13 = $ws.Range("O13").Value()
$P13 = $ws.Range("P13").Value()
$Q13 = $ws.Range("Q13").Value()

$D14 = $ws.Range("D14").Value()
$I14 = $ws.Range("I14").Value()
$J14 = $ws.Range("J14").Value()
$K14 = $ws.Range("K14").Value()
$L14 = $ws.Range("L14").Value()
$M14 = $ws.Range("M14").Value()
$N14 = $ws.Range("N14").Value()
$O14 = $ws.Range("O14").Value()
$P14 = $ws.Range("P14").Value()
$Q14 = $ws.Range("Q14").Value()

$D15 = $ws.Range("D15").Value()
$I15 = $ws.Range("I15").Value()
$J15 = $ws.Range("J15").Value()
$K15 = $ws.Range("K15").Value()
$L15 = $ws.Range("L15").Value()
$M15 = $ws.Range("M15").Value()
$N15 = $ws.Range("N15").Value()
$O15 = $ws.Range("O15").Value()
$P15 = $ws.Range("P15").Value()
$Q15 = $ws.Range("Q15").Value()

$D16 = $ws.Range("D16").Value()
$I16 = $ws.Range("I16").Value()
$J16 = $ws.Range("J16").Value()
$K16 = $ws.Range("K16").Value()
$L16 = $ws.Range("L16").Value()
$M16 = $ws.Range("M16").Value()
$N16 = $ws.Range("N16").Value()
$O16 = $ws.Range("O16").Value()
$P16 = $ws.Range("P16").Value()
$Q16 = $ws.Range("Q16").Value()

# Write permuted values: row[target] = original(row[source])
$ws.Range("D3").Value = $D9
$ws.Range("I3").Value = $I9
$ws.Range("J3").Value = $J9
$ws.Range("K3").Value = $K9
$ws.Range("L3").Value = $L9
$ws.Range("M3").Value = $M9
$ws.Range("N3").Value = $N9
$ws.Range("O3").Value = $O9
$ws.Range("P3").Value = $P9
$ws.Range("Q3").Value = $Q9

$ws.Range("D4").Value = $D13
$ws.Range("I4").Value = $I13
$ws.Range("J4").Value = $J13
$ws.Range("K4").Value = $K13
$ws.Range("L4").Value = $L13
$ws.Range("M4").Value = $M13
$ws.Range("N4").Value = $N13
$ws.Range("O4").Value = $O13
$ws.Range("P4").Value = $P13
$ws.Range("Q4").Value = $Q13

$ws.Range("D5").Value = $D14
$ws.Range("I5").Value = $I14
$ws.Range("J5").Value = $J14
$ws.Range("K5").Value = $K14
$ws.Range("L5").Value = $L14
$ws.Range("M5").Value = $M14
$ws.Range("N5").Value = $N14
$ws.Range("O5").Value = $O14
$ws.Range("P5").Value = $P14
$ws.Range("Q5").Value = $Q14

$ws.Range("D7").Value = $D4
$ws.Range("I7").Value = $I4
$ws.Range("J7").Value = $J4
$ws.Range("K7").Value = $K4
$ws.Range("L7").Value = $L4
$ws.Range("M7").Value = $M4
$ws.Range("N7").Value = $N4
$ws.Range("O7").Value = $O4
$ws.Range("P7").Value = $P4
$ws.Range("Q7").Value = $Q4

$ws.Range("D8").Value = $D10
$ws.Range("I8").Value = $I10
$ws.Range("J8").Value = $J10
$ws.Range("K8").Value = $K10
$ws.Range("L8").Value = $L10
$ws.Range("M8").Value = $M10
$ws.Range("N8").Value = $N10
$ws.Range("O8").Value = $O10
$ws.Range("P8").Value = $P10
$ws.Range("Q8").Value = $Q10

$ws.Range("D9").Value = $D11
$ws.Range("I9").Value = $I11
$ws.Range("J9").Value = $J11
$ws.Range("K9").Value = $K11
$ws.Range("L9").Value = $L11
$ws.Range("M9").Value = $M11
$ws.Range("N9").Value = $N11
$ws.Range("O9").Value = $O11
$ws.Range("P9").Value = $P11
$ws.Range("Q9").Value = $Q11

$ws.Range("D10").Value = $D16
$ws.Range("I10").Value = $I16
$ws.Range("J10").Value = $J16
$ws.Range("K10").Value = $K16
$ws.Range("L10").Value = $L16
$ws.Range("M10").Value = $M16
$ws.Range("N10").Value = $N16
$ws.Range("O10").Value = $O16
$ws.Range("P10").Value = $P16
$ws.Range("Q10").Value = $Q16

$ws.Range("D11").Value = $D12
$ws.Range("I11").Value = $I12
$ws.Range("J11").Value = $J12
$ws.Range("K11").Value = $K12
$ws.Range("L11").Value = $L12
$ws.Range("M11").Value = $M12
$ws.Range("N11").Value = $N12
$ws.Range("O11").Value = $O12
$ws.Range("P11").Value = $P12
$ws.Range("Q11").Value = $Q12

$ws.Range("D12").Value = $D5
$ws.Range("I12").Value = $I5
$ws.Range("J12").Value = $J5
$ws.Range("K12").Value = $K5
$ws.Range("L12").Value = $L5
$ws.Range("M12").Value = $M5
$ws.Range("N12").Value = $N5
$ws.Range("O12").Value = $O5
$ws.Range("P12").Value = $P5
$ws.Range("Q12").Value = $Q5

$ws.Range("D13").Value = $D15
$ws.Range("I13").Value = $I15
$ws.Range("J13").Value = $J15
$ws.Range("K13").Value = $K15
$ws.Range("L13").Value = $L15
$ws.Range("M13").Value = $M15
$ws.Range("N13").Value = $N15
$ws.Range("O13").Value = $O15
$ws.Range("P13").Value = $P15
$ws.Range("Q13").Value = $Q15

$ws.Range("D14").Value = $D3
$ws.Range("I14").Value = $I3
$ws.Range("J14").Value = $J3
$ws.Range("K14").Value = $K3
$ws.Range("L14").Value = $L3
$ws.Range("M14").Value = $M3
$ws.Range("N14").Value = $N3
$ws.Range("O14").Value = $O3
$ws.Range("P14").Value = $P3
$ws.Range("Q14").Value = $Q3

$ws.Range("D15").Value = $D7
$ws.Range("I15").Value = $I7
$ws.Range("J15").Value = $J7
$ws.Range("K15").Value = $K7
$ws.Range("L15").Value = $L7
$ws.Range("M15").Value = $M7
$ws.Range("N15").Value = $N7
$ws.Range("O15").Value = $O7
$ws.Range("P15").Value = $P7
$ws.Range("Q15").Value = $Q7

$ws.Range("D16").Value = $D8
$ws.Range("I16").Value = $I8
$ws.Range("J16").Value = $J8
$ws.Range("K16").Value = $K8
$ws.Range("L16").Value = $L8
$ws.Range("M16").Value = $M8
$ws.Range("N16").Value = $N8
$ws.Range("O16").Value = $O8
$ws.Range("P16").Value = $P8
$ws.Range("Q16").Value = $Q8
